# Generate Report for Handoff
# Updates the localization-status report with the latest handoff
# timestamps and sets the "Priority" column to "ht" for the files
# that were just handed off.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Rows corresponding to the files that were handed off in this run.
$rows = @(7, 8, 10, 11, 12, 14)

foreach ($r in $rows) {
    # "Latest HO Xliff Generate Date" on the Overview sheet.
    $wsOverview.Range("G$r").Value = "2016-08-29 00:22:06"

    # zh-cn sheet: new handoff datetime + priority flag.
    $wsZhCn.Range("H$r").Value = "2016-08-29 00:21:57"
    $wsZhCn.Range("E$r").Value = "ht"

    # de-de sheet: new handoff datetime + priority flag.
    $wsDeDe.Range("H$r").Value = "2016-08-29 00:22:06"
    $wsDeDe.Range("E$r").Value = "ht"
}
